$d = $word.ActiveDocument

# The new content must land right after the table and right before the
# final (bookmarked) paragraph that closes the document.
$count = $d.Paragraphs.Count
$anchorPara = $d.Paragraphs($count)
$insertionPoint = $anchorPara.Range.Duplicate
$insertionPoint.Collapse(1)   # wdCollapseStart

$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# Build: two blank paragraphs, then the descriptive paragraph (several runs,
# with proofing-error markers bracketing the two brand names), then one
# throw-away blank paragraph. InsertXML("...", "Before") merges the LAST
# paragraph of the inserted fragment into the anchor paragraph rather than
# giving it its own paragraph mark, so the trailing dummy paragraph absorbs
# that merge and is deleted afterwards, leaving the anchor paragraph (with
# the _GoBack bookmark) untouched.
$newParaXml = @"
<w:p $wNs><w:pPr><w:jc w:val="left"/></w:pPr><w:r><w:t xml:space="preserve">Podemos observar mediante el siguiente gráfico </w:t></w:r><w:r><w:t xml:space="preserve">cómo </w:t></w:r><w:r><w:t xml:space="preserve">consume </w:t></w:r><w:r><w:t xml:space="preserve">más memoria </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Adblock</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> Plus </w:t></w:r><w:r><w:t>respecto a</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>uBlock</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>.</w:t></w:r></w:p>
"@

$xml = "<w:p $wNs/><w:p $wNs/>" + $newParaXml + "<w:p $wNs/>"
$insertionPoint.InsertXML($xml, "Before")

# Drop the single paragraph mark belonging to the trailing dummy paragraph
# that InsertXML fused onto the bookmarked paragraph.
$bm = $d.Bookmarks("_GoBack")
$dummyMark = $d.Range($bm.Start - 1, $bm.Start)
$dummyMark.Delete()
